$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 55
$ws1.Range("F5").Value = 181
$ws1.Range("F6").Value = 9547
$ws1.Range("F9").Value = 1207
$ws1.Range("F10").Value = 1690
$ws1.Range("F11").Value = 156
$ws1.Range("F12").Value = 105
$ws1.Range("F14").Value = 270
$ws1.Range("F15").Value = 448
$ws1.Range("F16").Value = 96
$ws1.Range("F18").Value = 1314

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 55
$ws4.Range("F6").Value = 181
$ws4.Range("F7").Value = 9547
$ws4.Range("F10").Value = 1207
$ws4.Range("F11").Value = 1690
$ws4.Range("F12").Value = 156
$ws4.Range("F13").Value = 105
$ws4.Range("F15").Value = 270
$ws4.Range("F16").Value = 448
$ws4.Range("F17").Value = 96
$ws4.Range("F19").Value = 1314
